$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header for the DHCP/DNS MX-record preference value
$ws.Range("E1").Value2 = "mx-preference"

# New row 4: test3 (A record only, no mac-address)
$ws.Range("A4").Value2 = "test3"
$ws.Range("B4").Value2 = "int.domain.com"
$ws.Range("C4").Value2 = "192.168.1.2"

# New row 5: test4 (MX record, preference 5)
$ws.Range("A5").Value2 = "test4"
$ws.Range("B5").Value2 = "email.domain.com"
$ws.Range("C5").Value2 = "192.168.2.4"
$ws.Range("E5").Value2 = 5

# New row 6: test5 (MX record, preference 10)
$ws.Range("A6").Value2 = "test5"
$ws.Range("B6").Value2 = "email.domain.com"
$ws.Range("C6").Value2 = "192.168.2.5"
$ws.Range("E6").Value2 = 10

# Column width adjustments (col B widened for longer domain names, new col E added)
$ws.Columns.Item(2).ColumnWidth = 15.666666666666666
$ws.Columns.Item(5).ColumnWidth = 12.666666666666666

# Font charset tweak that accompanied the edit (Arial / ANSI charset)
$ws.Range("A1:E6").Font.Charset = 1

# Selection moves to E7 (next empty row) as in the authored file
$ws.Range("E7").Select()
